# P-TPI-Alexandre King.pptx
# 1) Mise a jour du journal : la date figee "21.03.2024" -> "23.05.2024"
#    (placeholder de date sur le masque de diapositive, les 11 mises en
#    page, le masque des documents et le masque des pages de notes).
# 2) Ajout d'une note pour la defense sur la diapositive 6 : le bloc
#    "Mise en place du MVC, Tailwind, uWamp, DB" devient "Parler du git"
#    et est redimensionne / repositionne en consequence.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "21.03.2024") {
                $shp.TextFrame.TextRange.Text = "23.05.2024"
            }
        }
    }
}

# -- Slide master --
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# -- Every slide layout ("mise en page") attached to the master --
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# -- Handout master --
Update-DatePlaceholder $p.HandoutMaster.Shapes

# -- Notes master --
Update-DatePlaceholder $p.NotesMaster.Shapes

# -- Slide 6 : "Deroulement" - replace the MVC/Tailwind/uWamp/DB note --
$slide6 = $p.Slides.Item(6)
$shape = $slide6.Shapes.Item(23)

$shape.TextFrame.TextRange.Text = "Parler du git"

# Reposition/resize the shrunk note box (values chosen so the stored EMU
# match 405534/4740763 offset and 1739752/223394 extent exactly).
$shape.Top = 373.2884521484375
$shape.Height = 17.590090
